$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, pushing the existing rows 12-17 down to 13-18.
$ws.Rows("12").Insert()

# Populate the new row 12 with the new weekly record.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 45089
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100104
$ws.Range("H12").Value = "Frutos de pepita"
$ws.Range("I12").Value = 100104005
$ws.Range("J12").Value = "Pera asiática"
$ws.Range("K12").Value = "Hosui"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 16000
$ws.Range("Q12").Value = "$/caja 18 kilos granel"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 889
$ws.Range("T12").Value = 18
